$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> C, D, E values (only cells with a numeric <v> in the sheet)
$data = @{
    3  = @(1, 1, 0.37)
    4  = @(0, 1, 0)
    5  = @(8, 5, 0.09)
    6  = @(3, 7, 0.02)
    7  = @(3, 4, 0.17)
    9  = @(41, 50, 0.02)
    11 = @(1, 3, 0.06)
    12 = @(2, 0, 0.14)
    13 = @(4, 0, 0.02)
    14 = @(2, 0, 0.14)
    17 = @(9, 12, 0.07)
    18 = @(1, 2, 0.18)
    19 = @(7, 7, 0.15)
    22 = @(1, 1, 0.37)
    23 = @(0, 0, 1)
    25 = @(6, 4, 0.13)
    26 = @(1, 0, 0.37)
    29 = @(0, 1, 0)
    32 = @(1, 0, 0.37)
    34 = @(6, 6, 0.16)
    35 = @(4, 1, 0.07)
    36 = @(8, 11, 0.07)
    37 = @(0, 0, 1)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
    $ws.Cells.Item($row, 5).Value = $vals[2]
}
